# Update Odd_* columns for rows 7-17 with refreshed FlashScore odds data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 7).Value = 2.35
$ws.Cells.Item(7, 8).Value = 3
$ws.Cells.Item(7, 9).Value = 3
$ws.Cells.Item(7, 11).Value = 1.91
$ws.Cells.Item(7, 12).Value = 3.7
$ws.Cells.Item(7, 13).Value = 1.04
$ws.Cells.Item(7, 14).Value = 7.1
$ws.Cells.Item(7, 15).Value = 1.42
$ws.Cells.Item(7, 16).Value = 2.47
$ws.Cells.Item(7, 17).Value = 2.18
$ws.Cells.Item(7, 18).Value = 1.53
$ws.Cells.Item(7, 20).Value = 2.25
$ws.Cells.Item(7, 21).Value = 1.9
$ws.Cells.Item(7, 22).Value = 1.72
$ws.Cells.Item(7, 23).Value = 6.6
$ws.Cells.Item(7, 24).Value = 10.5
$ws.Cells.Item(7, 29).Value = 7.3
$ws.Cells.Item(7, 31).Value = 16
$ws.Cells.Item(7, 33).Value = 900
$ws.Cells.Item(7, 34).Value = 7.7
$ws.Cells.Item(7, 35).Value = 14.5
$ws.Cells.Item(7, 36).Value = 11
$ws.Cells.Item(7, 38).Value = 30
$ws.Cells.Item(7, 39).Value = 45
$ws.Cells.Item(7, 40).Value = 4.05
$ws.Cells.Item(7, 42).Value = 25
$ws.Cells.Item(7, 46).Value = 2.22
$ws.Cells.Item(7, 47).Value = 7.7
$ws.Cells.Item(7, 48).Value = 90
$ws.Cells.Item(7, 49).Value = 4.7
$ws.Cells.Item(7, 50).Value = 17.5
$ws.Cells.Item(7, 51).Value = 29
$ws.Cells.Item(8, 7).Value = 2.07
$ws.Cells.Item(8, 8).Value = 3.15
$ws.Cells.Item(8, 9).Value = 3.45
$ws.Cells.Item(8, 10).Value = 2.62
$ws.Cells.Item(8, 11).Value = 2.07
$ws.Cells.Item(8, 12).Value = 3.9
$ws.Cells.Item(8, 14).Value = 7.2
$ws.Cells.Item(8, 15).Value = 1.33
$ws.Cells.Item(8, 16).Value = 2.8
$ws.Cells.Item(8, 17).Value = 1.98
$ws.Cells.Item(8, 18).Value = 1.65
$ws.Cells.Item(8, 22).Value = 1.85
$ws.Cells.Item(8, 23).Value = 7
$ws.Cells.Item(8, 26).Value = 19.5
$ws.Cells.Item(8, 27).Value = 17
$ws.Cells.Item(8, 28).Value = 29
$ws.Cells.Item(8, 29).Value = 8.5
$ws.Cells.Item(8, 30).Value = 6.1
$ws.Cells.Item(8, 36).Value = 12
$ws.Cells.Item(8, 41).Value = 10.25
$ws.Cells.Item(8, 42).Value = 18
$ws.Cells.Item(8, 43).Value = 40
$ws.Cells.Item(8, 44).Value = 70
$ws.Cells.Item(8, 46).Value = 2.55
$ws.Cells.Item(8, 48).Value = 60
$ws.Cells.Item(8, 50).Value = 19
$ws.Cells.Item(8, 51).Value = 25
$ws.Cells.Item(8, 52).Value = 100
$ws.Cells.Item(9, 7).Value = 1.47
$ws.Cells.Item(9, 8).Value = 4.25
$ws.Cells.Item(9, 9).Value = 5.1
$ws.Cells.Item(9, 10).Value = 1.95
$ws.Cells.Item(9, 11).Value = 2.42
$ws.Cells.Item(9, 12).Value = 4.85
$ws.Cells.Item(9, 13).Value = 1.01
$ws.Cells.Item(9, 14).Value = 14
$ws.Cells.Item(9, 15).Value = 1.12
$ws.Cells.Item(9, 16).Value = 4.65
$ws.Cells.Item(9, 17).Value = 1.48
$ws.Cells.Item(9, 18).Value = 2.33
$ws.Cells.Item(9, 19).Value = 1.25
$ws.Cells.Item(9, 20).Value = 3.6
$ws.Cells.Item(9, 21).Value = 1.62
$ws.Cells.Item(9, 22).Value = 2.24
$ws.Cells.Item(9, 23).Value = 7.8
$ws.Cells.Item(9, 24).Value = 7.2
$ws.Cells.Item(9, 25).Value = 7.1
$ws.Cells.Item(9, 26).Value = 9.25
$ws.Cells.Item(9, 27).Value = 9.25
$ws.Cells.Item(9, 28).Value = 16.5
$ws.Cells.Item(9, 29).Value = 16
$ws.Cells.Item(9, 30).Value = 7.7
$ws.Cells.Item(9, 31).Value = 12.5
$ws.Cells.Item(9, 32).Value = 40
$ws.Cells.Item(9, 33).Value = 250
$ws.Cells.Item(9, 34).Value = 15
$ws.Cells.Item(9, 35).Value = 27
$ws.Cells.Item(9, 36).Value = 13.5
$ws.Cells.Item(9, 37).Value = 65
$ws.Cells.Item(9, 38).Value = 35
$ws.Cells.Item(9, 39).Value = 32
$ws.Cells.Item(9, 40).Value = 3.55
$ws.Cells.Item(9, 41).Value = 6.8
$ws.Cells.Item(9, 42).Value = 14
$ws.Cells.Item(9, 43).Value = 19
$ws.Cells.Item(9, 44).Value = 40
$ws.Cells.Item(9, 45).Value = 150
$ws.Cells.Item(9, 46).Value = 3.35
$ws.Cells.Item(9, 47).Value = 7.3
$ws.Cells.Item(9, 48).Value = 55
$ws.Cells.Item(9, 49).Value = 7
$ws.Cells.Item(9, 50).Value = 26
$ws.Cells.Item(9, 51).Value = 28
$ws.Cells.Item(9, 52).Value = 150
$ws.Cells.Item(9, 53).Value = 150
$ws.Cells.Item(9, 54).Value = 300
$ws.Cells.Item(10, 7).Value = 5.1
$ws.Cells.Item(10, 8).Value = 3.95
$ws.Cells.Item(10, 9).Value = 1.52
$ws.Cells.Item(10, 10).Value = 5.1
$ws.Cells.Item(10, 11).Value = 2.25
$ws.Cells.Item(10, 12).Value = 2.05
$ws.Cells.Item(10, 13).Value = 1.02
$ws.Cells.Item(10, 14).Value = 8.300000000000001
$ws.Cells.Item(10, 15).Value = 1.23
$ws.Cells.Item(10, 16).Value = 3.42
$ws.Cells.Item(10, 17).Value = 1.72
$ws.Cells.Item(10, 18).Value = 1.9
$ws.Cells.Item(10, 19).Value = 1.36
$ws.Cells.Item(10, 20).Value = 2.99
$ws.Cells.Item(10, 21).Value = 1.85
$ws.Cells.Item(10, 22).Value = 1.91
$ws.Cells.Item(10, 23).Value = 11.75
$ws.Cells.Item(10, 24).Value = 24
$ws.Cells.Item(10, 25).Value = 14
$ws.Cells.Item(10, 26).Value = 70
$ws.Cells.Item(10, 27).Value = 40
$ws.Cells.Item(10, 28).Value = 40
$ws.Cells.Item(10, 29).Value = 11.25
$ws.Cells.Item(10, 30).Value = 6.8
$ws.Cells.Item(10, 31).Value = 14
$ws.Cells.Item(10, 32).Value = 60
$ws.Cells.Item(10, 33).Value = 400
$ws.Cells.Item(10, 34).Value = 6
$ws.Cells.Item(10, 35).Value = 6.1
$ws.Cells.Item(10, 36).Value = 7
$ws.Cells.Item(10, 37).Value = 9
$ws.Cells.Item(10, 38).Value = 10
$ws.Cells.Item(10, 39).Value = 21
$ws.Cells.Item(10, 40).Value = 6.7
$ws.Cells.Item(10, 41).Value = 29
$ws.Cells.Item(10, 42).Value = 35
$ws.Cells.Item(10, 43).Value = 175
$ws.Cells.Item(10, 44).Value = 200
$ws.Cells.Item(10, 45).Value = 450
$ws.Cells.Item(10, 46).Value = 2.85
$ws.Cells.Item(10, 47).Value = 7.9
$ws.Cells.Item(10, 48).Value = 75
$ws.Cells.Item(10, 49).Value = 3.35
$ws.Cells.Item(10, 50).Value = 7.2
$ws.Cells.Item(10, 51).Value = 17
$ws.Cells.Item(10, 52).Value = 22
$ws.Cells.Item(10, 53).Value = 55
$ws.Cells.Item(10, 54).Value = 250
$ws.Cells.Item(11, 7).Value = 18
$ws.Cells.Item(11, 8).Value = 8
$ws.Cells.Item(11, 9).Value = 1.08
$ws.Cells.Item(11, 10).Value = 13
$ws.Cells.Item(11, 11).Value = 3.3
$ws.Cells.Item(11, 12).Value = 1.34
$ws.Cells.Item(11, 15).Value = 1.06
$ws.Cells.Item(11, 16).Value = 8
$ws.Cells.Item(11, 17).Value = 1.25
$ws.Cells.Item(11, 18).Value = 3.5
$ws.Cells.Item(11, 19).Value = 1.14
$ws.Cells.Item(11, 20).Value = 5
$ws.Cells.Item(11, 21).Value = 2.55
$ws.Cells.Item(11, 22).Value = 1.49
$ws.Cells.Item(11, 23).Value = 50
$ws.Cells.Item(11, 24).Value = 175
$ws.Cells.Item(11, 25).Value = 60
$ws.Cells.Item(11, 27).Value = 300
$ws.Cells.Item(11, 28).Value = 175
$ws.Cells.Item(11, 29).Value = 22
$ws.Cells.Item(11, 30).Value = 18
$ws.Cells.Item(11, 31).Value = 37
$ws.Cells.Item(11, 32).Value = 175
$ws.Cells.Item(11, 34).Value = 9.5
$ws.Cells.Item(11, 35).Value = 6
$ws.Cells.Item(11, 36).Value = 11.25
$ws.Cells.Item(11, 37).Value = 5.3
$ws.Cells.Item(11, 38).Value = 10
$ws.Cells.Item(11, 39).Value = 35
$ws.Cells.Item(11, 40).Value = 18
$ws.Cells.Item(11, 41).Value = 120
$ws.Cells.Item(11, 42).Value = 90
$ws.Cells.Item(11, 46).Value = 4.4
$ws.Cells.Item(11, 47).Value = 12.5
$ws.Cells.Item(11, 48).Value = 120
$ws.Cells.Item(11, 49).Value = 3.1
$ws.Cells.Item(11, 50).Value = 4
$ws.Cells.Item(11, 51).Value = 16.5
$ws.Cells.Item(11, 52).Value = 7.4
$ws.Cells.Item(11, 53).Value = 30
$ws.Cells.Item(11, 54).Value = 250
$ws.Cells.Item(12, 7).Value = 4.75
$ws.Cells.Item(12, 9).Value = 1.7
$ws.Cells.Item(12, 10).Value = 5.5
$ws.Cells.Item(12, 12).Value = 2.3
$ws.Cells.Item(12, 15).Value = 1.36
$ws.Cells.Item(12, 16).Value = 3
$ws.Cells.Item(12, 17).Value = 2.15
$ws.Cells.Item(12, 18).Value = 1.67
$ws.Cells.Item(12, 21).Value = 2.1
$ws.Cells.Item(12, 22).Value = 1.67
$ws.Cells.Item(12, 23).Value = 12
$ws.Cells.Item(12, 28).Value = 51
$ws.Cells.Item(12, 29).Value = 8
$ws.Cells.Item(12, 41).Value = 29
$ws.Cells.Item(12, 42).Value = 41
$ws.Cells.Item(12, 44).Value = 151
$ws.Cells.Item(12, 45).Value = 500
$ws.Cells.Item(12, 50).Value = 9
$ws.Cells.Item(14, 15).Value = 1.29
$ws.Cells.Item(14, 16).Value = 3.5
$ws.Cells.Item(14, 17).Value = 1.93
$ws.Cells.Item(14, 18).Value = 1.93
$ws.Cells.Item(14, 23).Value = 6
$ws.Cells.Item(14, 25).Value = 8.5
$ws.Cells.Item(14, 29).Value = 10
$ws.Cells.Item(14, 34).Value = 17
$ws.Cells.Item(14, 35).Value = 41
$ws.Cells.Item(14, 36).Value = 23
$ws.Cells.Item(14, 53).Value = 201
$ws.Cells.Item(15, 14).Value = 9
$ws.Cells.Item(16, 12).Value = 3.75
$ws.Cells.Item(16, 17).Value = 2.4
$ws.Cells.Item(16, 18).Value = 1.53
$ws.Cells.Item(16, 19).Value = 1.53
$ws.Cells.Item(16, 20).Value = 2.38
$ws.Cells.Item(16, 31).Value = 17
$ws.Cells.Item(16, 33).Value = 1250
$ws.Cells.Item(16, 46).Value = 2.38
$ws.Cells.Item(17, 7).Value = 1.78
$ws.Cells.Item(17, 8).Value = 4.25
$ws.Cells.Item(17, 9).Value = 3.5
$ws.Cells.Item(17, 10).Value = 2.22
$ws.Cells.Item(17, 11).Value = 2.62
$ws.Cells.Item(17, 12).Value = 3.6
$ws.Cells.Item(17, 15).Value = 1.09
$ws.Cells.Item(17, 16).Value = 6.2
$ws.Cells.Item(17, 17).Value = 1.3
$ws.Cells.Item(17, 18).Value = 3.2
$ws.Cells.Item(17, 19).Value = 1.19
$ws.Cells.Item(17, 20).Value = 4.15
$ws.Cells.Item(17, 21).Value = 1.33
$ws.Cells.Item(17, 22).Value = 3.05
$ws.Cells.Item(17, 23).Value = 16
$ws.Cells.Item(17, 26).Value = 19
$ws.Cells.Item(17, 27).Value = 12
$ws.Cells.Item(17, 28).Value = 15
$ws.Cells.Item(17, 29).Value = 10.75
$ws.Cells.Item(17, 30).Value = 10.25
$ws.Cells.Item(17, 31).Value = 11.25
$ws.Cells.Item(17, 32).Value = 26
$ws.Cells.Item(17, 33).Value = 100
$ws.Cells.Item(17, 34).Value = 24
$ws.Cells.Item(17, 35).Value = 32
$ws.Cells.Item(17, 36).Value = 13.5
$ws.Cells.Item(17, 37).Value = 55
$ws.Cells.Item(17, 38).Value = 25
$ws.Cells.Item(17, 41).Value = 8.75
$ws.Cells.Item(17, 42).Value = 11.75
$ws.Cells.Item(17, 43).Value = 24
$ws.Cells.Item(17, 44).Value = 32
$ws.Cells.Item(17, 45).Value = 80
$ws.Cells.Item(17, 46).Value = 4.15
$ws.Cells.Item(17, 47).Value = 6.1
$ws.Cells.Item(17, 48).Value = 29
$ws.Cells.Item(17, 49).Value = 6.5
$ws.Cells.Item(17, 50).Value = 17
$ws.Cells.Item(17, 51).Value = 16
$ws.Cells.Item(17, 53).Value = 60
$ws.Cells.Item(17, 54).Value = 110
$ws.Cells.Item(17, 55).Value = 300
